$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM data changes NATMI LR-pair computations for Il34-Csf1r.
# Raw ligand values changed for sending cluster "ECs" (rows 2-4)
# Raw receptor values changed for target cluster "ECs" (rows 2, 5, 8)
# All derived specificity / edge-weight columns are recalculated downstream.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.024860333333333
$ws.Range("H2").Value = 6.074581
$ws.Range("I2").Value = 0.06209502815843435
$ws.Range("J2").Value = 0.06209502815843435
$ws.Range("M2").Value = 0.2649023333333334
$ws.Range("N2").Value = 0.7947070000000001
$ws.Range("O2").Value = 0.07626253349831308
$ws.Range("P2").Value = 0.07626253349831306
$ws.Range("Q2").Value = 0.5363902269741111
$ws.Range("R2").Value = 4.827512042767
$ws.Range("S2").Value = 0.004735524165011293
$ws.Range("T2").Value = 0.004735524165011292

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.024860333333333
$ws.Range("H3").Value = 6.074581
$ws.Range("I3").Value = 0.06209502815843435
$ws.Range("J3").Value = 0.06209502815843435
$ws.Range("O3").Value = 0.4194066525831247
$ws.Range("P3").Value = 0.4194066525831247
$ws.Range("Q3").Value = 2.949884028944444
$ws.Range("R3").Value = 26.5489562605
$ws.Range("S3").Value = 0.02604306790198382
$ws.Range("T3").Value = 0.02604306790198382

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.024860333333333
$ws.Range("H4").Value = 6.074581
$ws.Range("I4").Value = 0.06209502815843435
$ws.Range("J4").Value = 0.06209502815843435
$ws.Range("N4").Value = 5.255467
$ws.Range("O4").Value = 0.5043308139185624
$ws.Range("P4").Value = 0.5043308139185623
$ws.Range("Q4").Value = 3.547195553814111
$ws.Range("R4").Value = 31.924759984327
$ws.Range("S4").Value = 0.03131643609143924
$ws.Range("T4").Value = 0.03131643609143924

# Row 5
$ws.Range("I5").Value = 0.09055216274676732
$ws.Range("J5").Value = 0.09055216274676733
$ws.Range("M5").Value = 0.2649023333333334
$ws.Range("N5").Value = 0.7947070000000001
$ws.Range("O5").Value = 0.07626253349831308
$ws.Range("P5").Value = 0.07626253349831306
$ws.Range("Q5").Value = 0.7822090845148888
$ws.Range("R5").Value = 7.039881760634
$ws.Range("S5").Value = 0.00690573734482004
$ws.Range("T5").Value = 0.00690573734482004

# Row 6
$ws.Range("I6").Value = 0.09055216274676732
$ws.Range("J6").Value = 0.09055216274676733
$ws.Range("O6").Value = 0.4194066525831247
$ws.Range("P6").Value = 0.4194066525831247
$ws.Range("S6").Value = 0.037978179461784
$ws.Range("T6").Value = 0.03797817946178401

# Row 7
$ws.Range("I7").Value = 0.09055216274676732
$ws.Range("J7").Value = 0.09055216274676733
$ws.Range("N7").Value = 5.255467
$ws.Range("O7").Value = 0.5043308139185624
$ws.Range("P7").Value = 0.5043308139185623
$ws.Range("Q7").Value = 5.172817190194889
$ws.Range("R7").Value = 46.555354711754
$ws.Range("S7").Value = 0.04566824594016328
$ws.Range("T7").Value = 0.04566824594016328

# Row 8
$ws.Range("I8").Value = 0.8473528090947983
$ws.Range("J8").Value = 0.8473528090947984
$ws.Range("M8").Value = 0.2649023333333334
$ws.Range("N8").Value = 0.7947070000000001
$ws.Range("O8").Value = 0.07626253349831308
$ws.Range("P8").Value = 0.07626253349831306
$ws.Range("Q8").Value = 7.319616063911445
$ws.Range("R8").Value = 65.876544575203
$ws.Range("S8").Value = 0.06462127198848173
$ws.Range("T8").Value = 0.06462127198848173

# Row 9
$ws.Range("I9").Value = 0.8473528090947983
$ws.Range("J9").Value = 0.8473528090947984
$ws.Range("O9").Value = 0.4194066525831247
$ws.Range("P9").Value = 0.4194066525831247
$ws.Range("S9").Value = 0.3553854052193568
$ws.Range("T9").Value = 0.3553854052193569

# Row 10
$ws.Range("I10").Value = 0.8473528090947983
$ws.Range("J10").Value = 0.8473528090947984
$ws.Range("N10").Value = 5.255467
$ws.Range("O10").Value = 0.5043308139185624
$ws.Range("P10").Value = 0.5043308139185623
$ws.Range("Q10").Value = 48.40526216147144
$ws.Range("R10").Value = 435.647359453243
$ws.Range("S10").Value = 0.4273461318869598
$ws.Range("T10").Value = 0.4273461318869598
